# Update countries & provincias Spain
#
# This mirrors a refreshed COVID-19 data pull: the "Pais" sheet is kept
# sorted by column B ("Casos totales") descending, so when a country's
# total overtakes/falls behind its neighbour the two rows swap identity.
# Below, the numeric columns (B,C,D,E,G,H) get the new scraped values for
# every affected row, three adjacent row-pairs swap their country name in
# column A (because the row below overtook the row above after refresh),
# and the "last updated" banner in A1 is bumped to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- "last updated" banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 01:30"

# ---- country name swaps (re-sort after refresh) ----------------------------
# Japon overtook Bielorrusia
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Bielorrusia"
# Montenegro overtook Tunez
$ws.Range("A108").Value = "Montenegro"
$ws.Range("A109").Value = "Tunez"
# Trinidad yTobago overtook Bahamas
$ws.Range("A139").Value = "Trinidad yTobago"
$ws.Range("A140").Value = "Bahamas"

# ---- refreshed numeric data -------------------------------------------------
$ws.Range("B4").Value = 6633993
$ws.Range("C4").Value = 44346
$ws.Range("D4").Value = 3913714
$ws.Range("E4").Value = 2522928
$ws.Range("G4").Value = 1024
$ws.Range("H4").Value = 197351
$ws.Range("B6").Value = 4283978
$ws.Range("C6").Value = 44215
$ws.Range("D6").Value = 3530655
$ws.Range("E6").Value = 622849
$ws.Range("G6").Value = 899
$ws.Range("H6").Value = 130474
$ws.Range("B13").Value = 535705
$ws.Range("C13").Value = 11507
$ws.Range("E13").Value = 124436
$ws.Range("G13").Value = 241
$ws.Range("H13").Value = 11148
$ws.Range("B24").Value = 259725
$ws.Range("C24").Value = 1618
$ws.Range("E24").Value = 17002
$ws.Range("B29").Value = 135626
$ws.Range("C29").Value = 702
$ws.Range("D29").Value = 119674
$ws.Range("E29").Value = 6789
$ws.Range("B47").Value = 73901
$ws.Range("C47").Value = 680
$ws.Range("D47").Value = 65590
$ws.Range("E47").Value = 6899
$ws.Range("H47").Value = 1412
$ws.Range("B48").Value = 73784
$ws.Range("C48").Value = 193
$ws.Range("D48").Value = 72369
$ws.Range("E48").Value = 677
$ws.Range("G48").Value = 6
$ws.Range("H48").Value = 738
$ws.Range("B56").Value = 56017
$ws.Range("C56").Value = 188
$ws.Range("D56").Value = 43998
$ws.Range("E56").Value = 10943
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 1076
$ws.Range("B69").Value = 33860
$ws.Range("C69").Value = 1447
$ws.Range("E69").Value = 12270
$ws.Range("B88").Value = 13470
$ws.Range("C88").Value = 33
$ws.Range("D88").Value = 6731
$ws.Range("E88").Value = 5905
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 834
$ws.Range("B92").Value = 12003
$ws.Range("C92").Value = 137
$ws.Range("E92").Value = 1367
$ws.Range("B105").Value = 7479
$ws.Range("C105").Value = 26
$ws.Range("D105").Value = 5660
$ws.Range("E105").Value = 1595
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 224
$ws.Range("B106").Value = 7266
$ws.Range("C106").Value = 44
$ws.Range("D106").Value = 6786
$ws.Range("E106").Value = 319
$ws.Range("B108").Value = 6385
$ws.Range("C108").Value = 163
$ws.Range("D108").Value = 4460
$ws.Range("E108").Value = 1810
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 115
$ws.Range("B109").Value = 6259
$ws.Range("C109").Value = 377
$ws.Range("D109").Value = 1956
$ws.Range("E109").Value = 4200
$ws.Range("G109").Value = 4
$ws.Range("H109").Value = 103
$ws.Range("B119").Value = 4749
$ws.Range("C119").Value = 2
$ws.Range("E119").Value = 2862
$ws.Range("B138").Value = 2898
$ws.Range("C138").Value = 79
$ws.Range("D138").Value = 1451
$ws.Range("E138").Value = 1429
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = 18
$ws.Range("B139").Value = 2825
$ws.Range("C139").Value = 127
$ws.Range("D139").Value = 762
$ws.Range("E139").Value = 2013
$ws.Range("G139").Value = 7
$ws.Range("H139").Value = 50
$ws.Range("B140").Value = 2814
$ws.Range("C140").Value = 93
$ws.Range("D140").Value = 1220
$ws.Range("E140").Value = 1529
$ws.Range("H140").Value = 65
$ws.Range("B154").Value = 1773
$ws.Range("C154").Value = 14
$ws.Range("D154").Value = 1490
$ws.Range("E154").Value = 238
$ws.Range("B155").Value = 1763
$ws.Range("C155").Value = 13
$ws.Range("D155").Value = 1144
$ws.Range("E155").Value = 567
$ws.Range("G155").Value = 3
$ws.Range("H155").Value = 52
$ws.Range("D190").Value = 161
$ws.Range("E190").Value = 7
